# Apply the FonctionQualifiee -> SavoirFaire restructuring edit.
$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet updates -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-07-21T11:52:46+00:00"
$meta.Cells.Item(18, 2).Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/SavoirFaire"

# --- 2. Elements sheet updates --------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Row 3 used to describe "FonctionQualifiee.fonctionQualifiee"; it now documents
# the renamed "typeSavoirFaire" element (new short/definition/binding/base path).
$ws.Cells.Item(3, 1).Value = "FonctionQualifiee.typeSavoirFaire"
$ws.Cells.Item(3, 2).Value = "FonctionQualifiee.typeSavoirFaire"
$ws.Cells.Item(3, 12).Value = " Le type de savoir-faire (qualifications/autres attributions) d" + [char]0x00E9 + "signe par exemple:** une sp" + [char]0x00E9 + "cialit" + [char]0x00E9 + " ordinale (S);** une comp" + [char]0x00E9 + "tence (C);** etc."
$ws.Cells.Item(3, 13).Value = $ws.Cells.Item(3, 12).Value()
$ws.Cells.Item(3, 26).Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R04-TypeSavoirFaire?vs"
$ws.Cells.Item(3, 32).Value = "SavoirFaire.typeSavoirFaire"

# New rows 4-6: dateReconnaissance, dateAbandon, and the (re-added) plain
# fonctionQualifiee element. Copy row 3's formatting down first so the new
# rows keep the same cell style (border/wrap) as the rest of the table.
$ws.Range("A3:AJ3").Copy()
$ws.Range("A4:AJ6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$dateReconnaissanceDef = " Date " + [char]0x00E0 + " laquelle, l" + [char]0x2019 + "organisme donnant l" + [char]0x2019 + "autorisation d" + [char]0x2019 + "exercer une qualification a reconnu cette qualification ou date " + [char]0x00E0 + " laquelle l'attribution a " + [char]0x00E9 + "t" + [char]0x00E9 + " donn" + [char]0x00E9 + "e au professionnel."
$dateAbandonDef = " Date " + [char]0x00E0 + " laquelle le professionnel a d" + [char]0x00E9 + "clar" + [char]0x00E9 + " renoncer " + [char]0x00E0 + " l" + [char]0x2019 + "exercice d" + [char]0x2019 + "un savoir-faire ou date " + [char]0x00E0 + " laquelle il ne souhaite plus le faire appara" + [char]0x00EE + "tre."
$fonctionQualifieeDef = " Fonction qualifi" + [char]0x00E9 + "e."

# Row 4: FonctionQualifiee.dateReconnaissance
$ws.Cells.Item(4, 1).Value = "FonctionQualifiee.dateReconnaissance"
$ws.Cells.Item(4, 2).Value = "FonctionQualifiee.dateReconnaissance"
$ws.Cells.Item(4, 6).Value = "0"
$ws.Cells.Item(4, 7).Value = "1"
$ws.Cells.Item(4, 11).Value = "date`n"
$ws.Cells.Item(4, 12).Value = $dateReconnaissanceDef
$ws.Cells.Item(4, 13).Value = $dateReconnaissanceDef
$ws.Cells.Item(4, 32).Value = "SavoirFaire.dateReconnaissance"
$ws.Cells.Item(4, 33).Value = "0"
$ws.Cells.Item(4, 34).Value = "1"

# Row 5: FonctionQualifiee.dateAbandon
$ws.Cells.Item(5, 1).Value = "FonctionQualifiee.dateAbandon"
$ws.Cells.Item(5, 2).Value = "FonctionQualifiee.dateAbandon"
$ws.Cells.Item(5, 6).Value = "0"
$ws.Cells.Item(5, 7).Value = "1"
$ws.Cells.Item(5, 11).Value = "date`n"
$ws.Cells.Item(5, 12).Value = $dateAbandonDef
$ws.Cells.Item(5, 13).Value = $dateAbandonDef
$ws.Cells.Item(5, 32).Value = "SavoirFaire.dateAbandon"
$ws.Cells.Item(5, 33).Value = "0"
$ws.Cells.Item(5, 34).Value = "1"

# Row 6: FonctionQualifiee.fonctionQualifiee (moved back in, now as a plain Coding)
$ws.Cells.Item(6, 1).Value = "FonctionQualifiee.fonctionQualifiee"
$ws.Cells.Item(6, 2).Value = "FonctionQualifiee.fonctionQualifiee"
$ws.Cells.Item(6, 6).Value = "0"
$ws.Cells.Item(6, 7).Value = "1"
$ws.Cells.Item(6, 11).Value = "Coding`n"
$ws.Cells.Item(6, 12).Value = $fonctionQualifieeDef
$ws.Cells.Item(6, 13).Value = $fonctionQualifieeDef
$ws.Cells.Item(6, 24).Value = "preferred"
$ws.Cells.Item(6, 26).Value = "https://interop.esante.gouv.fr/ig/fhir/mos/ValueSet/fonctionQualifiee-vs"
$ws.Cells.Item(6, 32).Value = "FonctionQualifiee.fonctionQualifiee"
$ws.Cells.Item(6, 33).Value = "0"
$ws.Cells.Item(6, 34).Value = "1"

# --- 3. Column width bump for A/B (ID/Path) and Z (Binding Value Set), which
# now hold longer strings (bestFit autosize in the source workbook). -----------
$ws.Columns.Item(1).ColumnWidth = 31.43
$ws.Columns.Item(2).ColumnWidth = 31.43
$ws.Columns.Item(26).ColumnWidth = 68.97
